$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 110 currently only has the date (A110 = "16/11/2021").
# Fill in the Done / Progress / To-do columns for that day.
$ws.Range("B110").Value = "OpenMax : Revision,Image Video common"
$ws.Range("C110").Value = "Updating the notes"
$ws.Range("D110").Value = "Revision of C-DS-OS concepts"

# Additional "Done" notes continued below the same day's row, each on its
# own row in column B (matching the pattern used elsewhere in the sheet).
$ws.Range("B111").Value = "Low level android media APIs"
$ws.Range("B112").Value = "Media Player,Codecs"
$ws.Range("B113").Value = "Internal  discussion with teammates on yavta and run.sh file copmmands"

# Leave the selection on the last entered cell, as in the authored edit.
$ws.Range("B113").Select()
